$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1. "About" sheet updates
# ------------------------------------------------------------------
$about = $wb.Worksheets.Item("About")
$about.Range("B3").Value = "Calibration"
$about.Range("A12").Value = "We use the value from PNNL's GCAM model across vehicle technologies."
$about.Range("A13").ClearContents()
[void]$about.Range("A12").Select()

# ------------------------------------------------------------------
# 2. Insert new sheet "A54.tranSubsector_logit_revised" right after "About"
# ------------------------------------------------------------------
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $about)
$newSheet.Name = "A54.tranSubsector_logit_revised"

$newSheet.Range("A1").Value = "# File: A54.tranSubsector_logit_revised.csv"
$newSheet.Range("A2").Value = "# Title: Transportation default subsector logit exponents"
$newSheet.Range("A3").Value = "# Source: Documented in JIRA issue https://jira.pnnl.gov/jira/browse/JGCRI-358?src=confmacro."
$newSheet.Range("A4").Value = "# Units: Unitless"
$newSheet.Range("A5").Value = "# Column types: ccicc "
$newSheet.Range("A6").Value = "# ----------"

$newSheet.Range("A7").Value = "supplysector"
$newSheet.Range("B7").Value = "tranSubsector"
$newSheet.Range("C7").Value = "logit.exponent"
$newSheet.Range("D7").Value = "logit.factor"
$newSheet.Range("E7").Value = "logit.type"

$newSheet.Range("A8").Value = "trn_aviation_intl"
$newSheet.Range("B8").Value = "International Aviation"
$newSheet.Range("C8").Value = -6
$newSheet.Range("E8").Value = "absolute-cost-logit"

$newSheet.Range("A9").Value = "trn_freight"
$newSheet.Range("B9").Value = "Domestic Ship"
$newSheet.Range("C9").Value = -6
$newSheet.Range("E9").Value = "absolute-cost-logit"

$newSheet.Range("A10").Value = "trn_freight"
$newSheet.Range("B10").Value = "Freight Rail"
$newSheet.Range("C10").Value = -1
$newSheet.Range("D10").Value = "Fuel types and efficiency levels"
$newSheet.Range("E10").Value = "absolute-cost-logit"

$newSheet.Range("A11").Value = "trn_freight_road"
$newSheet.Range("B11").Value = "Light truck"
$newSheet.Range("C11").Value = -8
$newSheet.Range("D11").Value = "Fuel types"
$newSheet.Range("E11").Value = "absolute-cost-logit"

$newSheet.Range("A12").Value = "trn_freight_road"
$newSheet.Range("B12").Value = "Medium truck"
$newSheet.Range("C12").Value = -8
$newSheet.Range("D12").Value = "Fuel types"
$newSheet.Range("E12").Value = "absolute-cost-logit"

$newSheet.Range("A13").Value = "trn_freight_road"
$newSheet.Range("B13").Value = "Heavy truck"
$newSheet.Range("C13").Value = -8
$newSheet.Range("D13").Value = "Fuel types"
$newSheet.Range("E13").Value = "absolute-cost-logit"

$newSheet.Range("A14").Value = "trn_pass"
$newSheet.Range("B14").Value = "Cycle"
$newSheet.Range("C14").Value = -6
$newSheet.Range("E14").Value = "absolute-cost-logit"

$newSheet.Range("A15").Value = "trn_pass"
$newSheet.Range("B15").Value = "Domestic Aviation"
$newSheet.Range("C15").Value = -6
$newSheet.Range("E15").Value = "absolute-cost-logit"

$newSheet.Range("A16").Value = "trn_pass"
$newSheet.Range("B16").Value = "HSR"
$newSheet.Range("C16").Value = -6
$newSheet.Range("E16").Value = "absolute-cost-logit"

$newSheet.Range("A17").Value = "trn_pass"
$newSheet.Range("B17").Value = "Passenger Rail"
$newSheet.Range("C17").Value = -1
$newSheet.Range("D17").Value = "Fuel types and efficiency levels"
$newSheet.Range("E17").Value = "absolute-cost-logit"

$newSheet.Range("A18").Value = "trn_pass"
$newSheet.Range("B18").Value = "Walk"
$newSheet.Range("C18").Value = -6
$newSheet.Range("E18").Value = "absolute-cost-logit"

$newSheet.Range("A19").Value = "trn_pass_road"
$newSheet.Range("B19").Value = "Bus"
$newSheet.Range("C19").Value = -3
$newSheet.Range("D19").Value = "Fuel types"
$newSheet.Range("E19").Value = "absolute-cost-logit"

$newSheet.Range("A20").Value = "trn_pass_road_LDV"
$newSheet.Range("B20").Value = "2W and 3W"
$newSheet.Range("C20").Value = -8
$newSheet.Range("D20").Value = "Fuel types"
$newSheet.Range("E20").Value = "absolute-cost-logit"

$newSheet.Range("A21").Value = "trn_pass_road_LDV_4W"
$newSheet.Range("B21").Value = "Car"
$newSheet.Range("C21").Value = -8
$newSheet.Range("D21").Value = "Fuel types and ICE efficiency levels"
$newSheet.Range("E21").Value = "absolute-cost-logit"

$newSheet.Range("A22").Value = "trn_pass_road_LDV_4W"
$newSheet.Range("B22").Value = "Large Car and Truck"
$newSheet.Range("C22").Value = -8
$newSheet.Range("D22").Value = "Fuel types and ICE efficiency levels"
$newSheet.Range("E22").Value = "absolute-cost-logit"

$newSheet.Range("A23").Value = "trn_pass_road_LDV_4W"
$newSheet.Range("B23").Value = "Mini Car"
$newSheet.Range("C23").Value = -8
$newSheet.Range("D23").Value = "Fuel types and ICE efficiency levels"
$newSheet.Range("E23").Value = "absolute-cost-logit"

$newSheet.Range("A24").Value = "trn_shipping_intl"
$newSheet.Range("B24").Value = "International Ship"
$newSheet.Range("C24").Value = -6
$newSheet.Range("E24").Value = "absolute-cost-logit"

$newSheet.Range("A25").Value = "# Passthrough tranSubsectors are listed belowabsolute-cost-logit"

$newSheet.Range("A26").Value = "trn_pass"
$newSheet.Range("B26").Value = "road"
$newSheet.Range("C26").Value = -6
$newSheet.Range("E26").Value = "absolute-cost-logit"

$newSheet.Range("A27").Value = "trn_pass_road"
$newSheet.Range("B27").Value = "LDV"
$newSheet.Range("C27").Value = -6
$newSheet.Range("E27").Value = "absolute-cost-logit"

$newSheet.Range("A28").Value = "trn_pass_road"
$newSheet.Range("B28").Value = "bus"
$newSheet.Range("C28").Value = -6
$newSheet.Range("E28").Value = "absolute-cost-logit"

$newSheet.Range("A29").Value = "trn_pass_road_LDV"
$newSheet.Range("B29").Value = "4W"
$newSheet.Range("C29").Value = -6
$newSheet.Range("E29").Value = "absolute-cost-logit"

$newSheet.Range("A30").Value = "trn_pass_road_LDV"
$newSheet.Range("B30").Value = "2W"
$newSheet.Range("C30").Value = -6
$newSheet.Range("E30").Value = "absolute-cost-logit"

$newSheet.Range("A31").Value = "trn_freight"
$newSheet.Range("B31").Value = "road"
$newSheet.Range("C31").Value = -6
$newSheet.Range("E31").Value = "absolute-cost-logit"

$newSheet.Columns.Item(1).ColumnWidth = 21.592447916666668
$newSheet.Columns.Item(2).ColumnWidth = 33.307291666666664
$newSheet.Columns.Item(3).ColumnWidth = 13.451822916666666

[void]$newSheet.Range("C56").Select()

# ------------------------------------------------------------------
# 3. "TTLE" sheet: link logit exponent values to the new sheet via formulas
# ------------------------------------------------------------------
$ttle = $wb.Worksheets.Item("TTLE")

$ttle.Range("B2").Formula = "='A54.tranSubsector_logit_revised'!C21"
$ttle.Range("C2").Formula = "='A54.tranSubsector_logit_revised'!C11"

$ttle.Range("B3").Formula = "='A54.tranSubsector_logit_revised'!C19"
$ttle.Range("C3").Formula = "='A54.tranSubsector_logit_revised'!C13"

$ttle.Range("B4").Formula = "='A54.tranSubsector_logit_revised'!C15"
$ttle.Range("C4").Formula = "='A54.tranSubsector_logit_revised'!C8"

$ttle.Range("B5").Formula = "='A54.tranSubsector_logit_revised'!C17"
$ttle.Range("C5").Formula = "='A54.tranSubsector_logit_revised'!C10"

$ttle.Range("B6").Formula = "='A54.tranSubsector_logit_revised'!C24"
$ttle.Range("C6").Formula = "='A54.tranSubsector_logit_revised'!C9"

$ttle.Range("B7").Formula = "='A54.tranSubsector_logit_revised'!C20"
$ttle.Range("C7").Formula = "='A54.tranSubsector_logit_revised'!C20"

[void]$ttle.Range("C8").Select()

# ------------------------------------------------------------------
# 4. Restore "About" as the active/selected sheet (matches original file
#    state where the About tab was the one left selected)
# ------------------------------------------------------------------
[void]$about.Range("A12").Select()
